# Update reference (previous_count) and change values for two agencies
# to reflect the refreshed comparison data ("changed reference file for
# archive to reflect changes since early am Dec 1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Department of Health and Human Services
$ws.Range("C10").Value = 44
$ws.Range("D10").Value = 0

# Row 26: Intelligence Community
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 0
